# Updated cryptos list on Sat Feb 18 18:29:53 UTC 2023 with GitHub Actions
#
# Refreshes the Price / Volume(1h) snapshot columns for every coin row.
# Row 33 and row 34 also swap identities this cycle (Hedera overtakes
# ImmutableX in the ranking), so their Coin name, Link and data columns
# are rewritten together.
#
# Cell values are written with a leading apostrophe (text qualifier) and
# the style is reset to "Normal" afterwards so that numeric-looking price
# strings (e.g. "0.9989") stay stored as text, matching the source feed's
# inline-string cells, instead of being auto-coerced to numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "24.760.12"
Set-TextValue 2 5 "  +1.98%  "
Set-TextValue 3 4 "1.700.76"
Set-TextValue 3 5 "  +0.60%  "
Set-TextValue 4 4 "0.9989"
Set-TextValue 4 5 "  -0.09%  "
Set-TextValue 5 4 "316.03"
Set-TextValue 5 5 "  +1.53%  "
Set-TextValue 6 4 "0.9990"
Set-TextValue 6 5 "  +0.04%  "
Set-TextValue 7 4 "0.3949"
Set-TextValue 7 5 "  +0.96%  "
Set-TextValue 8 4 "0.4059"
Set-TextValue 8 5 "  +0.82%  "
Set-TextValue 9 4 "1.494"
Set-TextValue 9 5 "  +2.14%  "
Set-TextValue 10 4 "0.9985"
Set-TextValue 10 5 "  -0.16%  "
Set-TextValue 11 4 "53.12"
Set-TextValue 11 5 "  -1.34%  "
Set-TextValue 12 4 "0.08886"
Set-TextValue 12 5 "  +1.95%  "
Set-TextValue 13 4 "7.241"
Set-TextValue 13 5 "  -0.57%  "
Set-TextValue 14 4 "23.71"
Set-TextValue 14 5 "  +3.21%  "
Set-TextValue 15 4 "8.103"
Set-TextValue 15 5 "  +9.03%  "
Set-TextValue 16 4 "0.00001324"
Set-TextValue 16 5 "  +0.67%  "
Set-TextValue 17 4 "1.697.45"
Set-TextValue 17 5 "  +1.09%  "
Set-TextValue 18 4 "99.95"
Set-TextValue 18 5 "  +0.02%  "
Set-TextValue 19 4 "0.07003"
Set-TextValue 19 5 "  -0.20%  "
Set-TextValue 20 4 "19.69"
Set-TextValue 20 5 "  +2.09%  "
Set-TextValue 21 4 "7.066"
Set-TextValue 21 5 "  +5.60%  "
Set-TextValue 22 4 "0.9991"
Set-TextValue 22 5 "  +0.19%  "
Set-TextValue 23 4 "14.40"
Set-TextValue 23 5 "  +1.95%  "
Set-TextValue 24 4 "24.742.05"
Set-TextValue 24 5 "  +1.88%  "
Set-TextValue 25 4 "3.295"
Set-TextValue 25 5 "  +11.01%  "
Set-TextValue 26 5 "  +1.31%  "
Set-TextValue 27 5 "  +1.45%  "
Set-TextValue 28 4 "163.47"
Set-TextValue 28 5 "  +2.52%  "
Set-TextValue 29 4 "136.42"
Set-TextValue 29 5 "  +2.61%  "
Set-TextValue 30 5 "  +1.52%  "
Set-TextValue 31 4 "7.466"
Set-TextValue 31 5 "  +1.09%  "
Set-TextValue 32 4 "1.884.71"
Set-TextValue 32 5 "  +1.08%  "
Set-TextValue 33 2 "Hedera"
Set-TextValue 33 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 33 4 "0.08620"
Set-TextValue 33 5 "  -0.54%  "
Set-TextValue 34 2 "ImmutableX"
Set-TextValue 34 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 34 4 "1.065"
Set-TextValue 34 5 "  -0.35%  "
Set-TextValue 35 4 "7.171"
Set-TextValue 35 5 "  -3.74%  "
Set-TextValue 36 4 "11.65"
Set-TextValue 36 5 "  +6.63%  "
Set-TextValue 37 4 "0.2757"
Set-TextValue 37 5 "  +2.28%  "
Set-TextValue 38 4 "1.933"
Set-TextValue 38 5 "  -1.15%  "
Set-TextValue 39 4 "14.58"
Set-TextValue 39 5 "  +0.11%  "
Set-TextValue 40 4 "0.09218"
Set-TextValue 40 5 "  +3.48%  "
Set-TextValue 41 4 "0.02743"
Set-TextValue 41 5 "  +2.40%  "
Set-TextValue 42 4 "1.473"
Set-TextValue 42 5 "  +1.37%  "
Set-TextValue 43 5 "  +1.75%  "
Set-TextValue 44 4 "16.10"
Set-TextValue 44 5 "  +4.68%  "
Set-TextValue 45 4 "0.7215"
Set-TextValue 45 5 "  +1.47%  "
Set-TextValue 46 4 "2.583"
Set-TextValue 46 5 "  +6.52%  "
Set-TextValue 47 4 "4.212"
Set-TextValue 47 5 "  +1.75%  "
Set-TextValue 48 4 "0.9989"
Set-TextValue 48 5 "  +0.09%  "
Set-TextValue 49 4 "1.334"
Set-TextValue 49 5 "  +5.08%  "
Set-TextValue 50 4 "140.04"
Set-TextValue 50 5 "  +0.24%  "
Set-TextValue 51 4 "0.08009"
Set-TextValue 51 5 "  +1.07%  "
